# Update "想去人数" (F column) counts on two sheets that share the same
# underlying data: "展览" and "全部类型".
# Row 6: F6 265 -> 266
# Row 7: F7 6473 -> 6474

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 266
    $ws.Range("F7").Value = 6474
}
